$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The old sheet had a two-row header (row1 + row2) on top of 10 data rows
# (rows 3-12). The new layout uses a single header row (row1) with five
# extra leading columns (idx, idx2, Name, Date Start, Date End) followed by
# the original measurement columns, and the 10 data rows shift up to rows
# 2-11. Deleting the old second header row achieves that shift while
# keeping the data + styles of the remaining rows intact.
$ws.Rows(2).Delete()

# Clear any left-over formatting on the new leading header cells (A1:E1)
# so they pick up the default/general style, matching the new columns.
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Rewrite the measurement header row (previously split across two rows)
# as a single row of labels.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# These headers use the small Arial 9pt font (same font used elsewhere in
# the sheet for labels / text cells).
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Match the saved selection state.
$null = $ws.Range("A2:K2").Select()
